$d = $word.ActiveDocument

# The document ends with two empty paragraphs after the "pickle" / reddit
# hyperlink paragraph. We need to insert two new paragraphs ("Sail Boat"
# and the pinimg URL) right after the reddit hyperlink paragraph, i.e.
# right before those two trailing empty paragraphs.

$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "reddit\.com/media") {
        $anchor = $p
    }
}

$r = $anchor.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$idx = $anchor.Index
$sailPara = $d.Paragraphs.Item($idx + 1)
$sailPara.Range.InsertAfter("Sail Boat")
$sailPara.Range.InsertParagraphAfter()

$urlPara = $d.Paragraphs.Item($idx + 2)
$urlPara.Range.InsertAfter("https://i.pinimg.com/736x/f4/36/ec/f436eccf30ca4093c0f26ed388dcd34e.jpg")
